$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update BASE AMOUNT (F), INITIAL AMOUNT (G) and TOTAL (H) for rows 2-10
$ws.Range("F2:F10").Value = 600

$ws.Range("G2").Value = 360
$ws.Range("G3").Value = 360
$ws.Range("G4").Value = 360

$ws.Range("H2").Value = 960
$ws.Range("H3").Value = 960
$ws.Range("H4").Value = 960

$ws.Range("H5:H10").Value = 600
